$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232") entirely - remaining rows shift up by one.
$ws.Rows(26).Delete()

# Delete the row that is now at 27 (originally "SC 92") - remaining rows shift up again.
$ws.Rows(27).Delete()

# Clear / fill in the column F values that moved (missing-data swap).
$ws.Range("F5").Value = ""
$ws.Range("F8").Value = 17.05
$ws.Range("F12").Value = ""
$ws.Range("F14").Value = 17.76
$ws.Range("F18").Value = ""

# Row 26 is now "SC 5" - fill in column B.
$ws.Range("B26").Value = -20.2

# Row 27 is now "SC 101" - clear column B.
$ws.Range("B27").Value = ""

# Row 33 is now "SC 232" - fill in column D.
$ws.Range("D33").Value = -14.1
